$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    'karen Regina Matos Cunha',
    'Evelyn Ribeira Soares',
    'Selma Alves de Alencar',
    'Monica Ester Miguel',
    'Aline de Carvalho Almeida',
    'Dariany Mickelly de Araujo Silva',
    'Jessica Soares de Oliveira',
    'Daiane Moreno Lima',
    'Giovanna Santos ancelmo',
    'Bruna Camargo Garnes',
    'Rafaela da Silva Amorim Santos',
    'Beatriz da Silva Oliveira',
    'Maiara Parisi lobo',
    'Daiane da Silva Mendes Nascimento',
    'Fabiana Chaves de Paula',
    'William Cesar Silva Alves',
    'Ana Paula Souza',
    'Alessandra Gomes do Nascimento',
    'Renata Pereira',
    'Guiomar Marim',
    'Lucilene da Silva Maciel',
    'Priscila da Costa Reis',
    'Thayna Conceição da Cunha',
    'Alessandra Floriano Silva',
    'Milena de Assis Marin',
    'Angelita arecida da Silva'
)

$quoted = $names | ForEach-Object { "'$_'" }
$namesList = "[" + ($quoted -join ", ") + "]"

$prefix = "`\Olá, pessoal! Sou um robô super simpático aqui para ajudar com a escala de folgas. Por favor, preencham o dia que desejam folgar, seguindo o exemplo abaixo:`n`nExemplo: Juriscreide Nascimento Bezerra - dia 28`n`nAgora é a vez de vocês:`n"

$ws.Range("A2").Value = $prefix + $namesList
